# Sync attendance_reports: normalize "Recorded By" (column G) so that a
# leading "System" entry is moved to the end of the comma-separated list
# instead of the front (e.g. "System, foo@bar.com" -> "foo@bar.com, System").
# Note: this runtime's -ceq / -cmatch operators do NOT behave case-sensitively,
# so an explicit char-by-char comparison is used to exactly match "System"
# (capital S) and leave a lowercase "system" token untouched.

function Test-ExactSystemToken($token) {
    if ($token.Length -ne 6) { return $false }
    $target = "System"
    for ($i = 0; $i -lt 6; $i++) {
        if ([int][char]$token[$i] -ne [int][char]$target[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1 -and (Test-ExactSystemToken $parts[0])) {
            $rest = $parts[1..($parts.Count - 1)]
            $newVal = ($rest + $parts[0]) -join ", "
            $cell.Value2 = $newVal
        }
    }
}
